$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 756, shifting existing rows 756-859 down to 757-860
$ws.Rows.Item(756).Insert()

# Populate the newly inserted row 756 with data
$ws.Range("A756").Value = 9
$ws.Range("B756").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C756").Value = "Metropolitana"
$ws.Range("D756").Value = 45212
$ws.Range("E756").Value = 13
$ws.Range("F756").Value = 100112031
$ws.Range("G756").Value = "Poroto verde"
$ws.Range("H756").Value = "Magnum"
$ws.Range("I756").Value = "Primera"
$ws.Range("J756").Value = 70
$ws.Range("K756").Value = 21000
$ws.Range("L756").Value = 23000
$ws.Range("M756").Value = 22000
$ws.Range("N756").Value = "`$/malla 25 kilos"
$ws.Range("O756").Value = "Perú"
$ws.Range("P756").Value = 880
$ws.Range("Q756").Value = 25
$ws.Range("R756").Value = "Hortaliza"
